$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Device")

# Delete the "Notes:" row, shifting everything up by one row.
$ws.Rows.Item(2).Delete()

# Rename the header cells from "Device type"/"Device parameters" to "Type"/"Parameters".
$ws.Range("B3").Value = "Type"
$ws.Range("C3").Value = "Parameters"

# Update the descriptive text at the top of the sheet.
$ws.Range("A1").Value = "This sheet summarizes the apparatuses connected to buses."

# Select row 2 (matches the author's final selection state in the saved file).
$ws.Rows.Item(2).Select()

# Rename the sheet tab itself.
$ws.Name = "Apparatus"
